$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A5:E9")

# Force the new cells to be treated as text (matching the rest of the
# sheet, which stores every value -- including numeric-looking ones --
# as text) instead of Excel's usual "smart" numeric auto-detection.
$range.NumberFormat = "@"

$data = @(
    @("10G108024", "Bag Poly - Turkey 10x8x24", "1", "33.98", "33.98"),
    @("711603", "NABC Bathroom Cleaner", "1", "35.35", "35.35"),
    @("315904", "Sani-T-10 Sanitizer", "1", "114.86", "114.86"),
    @("K8", "Wrap Poly 8x10.75", "4", "51.94", "207.76"),
    @("LKC1624F", "Parfait Lid", "1", "47.53", "47.53")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

# Drop the temporary "@" number format so the new rows keep the same
# (default) cell style as the rest of the sheet, now that the values
# are locked in as text.
$range.ClearFormats()
